$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushing old row4 "Mary Doe" and row5 "Steven Doe" down)
$ws.Rows.Item(4).Insert()

# Header row - add "notSupported" label in E1
$ws.Range("E1").Value = "notSupported"

# New row 4 content: id=3, name="42", country/language empty, notSupported="Extra"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "42"
$ws.Range("E4").Value = "Extra"

# Update ids for the shifted rows (5 and 6)
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Column E width
$ws.Columns.Item(5).ColumnWidth = 27.125
